$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the beginning so the existing Method / Average
# Significant Percentage columns (A,B) shift right to (C,D).
$ws.Range("A:B").Insert()

# Copy the header formatting (bold font + border) from the shifted header
# cells onto the two new header cells.
$ws.Range("C1:D1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# New header labels
$ws.Range("A1").Value = "Network"
$ws.Range("B1").Value = "Alpha"

# New Network / Alpha columns for each data row
$ws.Range("A2").Value = "HumanNet"
$ws.Range("A3").Value = "HumanNet"
$ws.Range("A4").Value = "HumanNet"
$ws.Range("A5").Value = "HumanNet"

$ws.Range("B2").Value = 0.2
$ws.Range("B3").Value = 0.2
$ws.Range("B4").Value = 0.2
$ws.Range("B5").Value = 0.2

# Updated Average Significant Percentage values (now column D)
$ws.Range("D2").Value = 0.02525252525252525
$ws.Range("D3").Value = 0.02270663033605813
$ws.Range("D4").Value = 0.04668620813034981
$ws.Range("D5").Value = 0.02397957779429169
